$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.369.28'
$ws.Range("E2").Value = '  +0.74%  '
$ws.Range("D3").Value = '2.643.05'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.56'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.83'
$ws.Range("E6").Value = '  +0.80%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.547'
$ws.Range("E8").Value = '  -0.52%  '
$ws.Range("D9").Value = '2.641.42'
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("E10").Value = '  +7.22%  '
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.27'
$ws.Range("E12").Value = '  +1.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +2.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.09'
$ws.Range("E14").Value = '  +1.69%  '
$ws.Range("E15").Value = '  +2.46%  '
$ws.Range("D16").Value = '3.122.78'
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").Value = '68.237.92'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '2.640.25'
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '363.83'
$ws.Range("E20").Value = '  -2.39%  '
$ws.Range("E21").Value = '  +0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.38'
$ws.Range("E22").Value = '  +2.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.89'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.50'
$ws.Range("E25").Value = '  +3.27%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.81'
$ws.Range("E27").Value = '  -1.45%  '
$ws.Range("E28").Value = '  +1.80%  '
$ws.Range("D29").Value = '2.777.04'
$ws.Range("E29").Value = '  +0.61%  '
$ws.Range("E30").Value = '  +0.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '572.20'
$ws.Range("E31").Value = '  -0.29%  '
$ws.Range("E32").Value = '  +3.23%  '
$ws.Range("E33").Value = '  +1.77%  '
$ws.Range("E34").Value = '  +1.47%  '
$ws.Range("E35").Value = '  +3.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  +5.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.06'
$ws.Range("E38").Value = '  +1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.37'
$ws.Range("E39").Value = '  +1.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.375'
$ws.Range("E40").Value = '  +1.50%  '
$ws.Range("E41").Value = '  -0.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.38'
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("E43").Value = '  +3.67%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.67'
$ws.Range("E44").Value = '  +1.90%  '
$ws.Range("E45").Value = '  +3.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.69'
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("E47").Value = '  +0.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '156.16'
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.76'
$ws.Range("E49").Value = '  +1.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.72'
$ws.Range("E50").Value = '  +0.96%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0789'
$ws.Range("E51").Value = '  +1.48%  '
